{"js": "// Replace every arithmetic-expression cell in the worksheet table with its\n// new expression. The mapping below is old-expression -> new-expression,\n// taken from the authoritative diff; using a lookup keyed on the existing\n// cell text (rather than blind positional indices) makes the edit robust to\n// the cells being visited in any order.\nconst replacements = {\n  \"44+39=\": \"39+45=\",\n  \"65+10=\": \"9-5=\",\n  \"56-27=\": \"13-0=\",\n  \"77+0=\": \"33+61=\",\n  \"13+9=\": \"17+19=\",\n  \"88-63=\": \"76-35=\",\n  \"59-6=\": \"78-9=\",\n  \"63-36=\": \"12+73=\",\n  \"21+54=\": \"51-17=\",\n  \"26+73=\": \"86-66=\",\n  \"52-37=\": \"6+49=\",\n  \"37+19=\": \"16-1=\",\n  \"84-15=\": \"49-44=\",\n  \"83-68=\": \"42+35=\",\n  \"5+33=\": \"6+75=\",\n  \"97-5=\": \"92-24=\",\n  \"11+44=\": \"20-8=\",\n  \"30-26=\": \"32+37=\",\n  \"66-52=\": \"8+17=\",\n  \"1+67=\": \"10+20=\",\n  \"22-17=\": \"20-10=\",\n  \"37+20=\": \"51-42=\",\n  \"4+88=\": \"2+61=\",\n  \"26+42=\": \"54+12=\",\n  \"57-47=\": \"48-10=\",\n  \"14-7=\": \"20+49=\",\n  \"11-0=\": \"11+63=\",\n  \"82-5=\": \"12+41=\",\n  \"91-66=\": \"64+19=\",\n  \"66-44=\": \"97-27=\",\n  \"36+31=\": \"98-24=\",\n  \"83-67=\": \"78-31=\",\n  \"66-4=\": \"91-8=\",\n  \"65+11=\": \"1+92=\",\n  \"37+32=\": \"64-34=\",\n  \"51-5=\": \"14+79=\",\n  \"91-52=\": \"96-82=\",\n  \"33+37=\": \"61-36=\",\n  \"57-7=\": \"67-23=\",\n  \"1+61=\": \"13+61=\",\n  \"32+60=\": \"23+45=\",\n  \"39+43=\": \"83+2=\",\n  \"22+58=\": \"46-4=\",\n  \"94-6=\": \"84-43=\",\n  \"37+25=\": \"4+66=\",\n  \"64-20=\": \"47-11=\",\n  \"32+4=\": \"49+8=\",\n  \"13+69=\": \"66-49=\",\n  \"45-8=\": \"88-22=\",\n  \"72-52=\": \"63-43=\",\n  \"29+8=\": \"56+5=\",\n  \"38+31=\": \"25-7=\",\n  \"75-25=\": \"14+33=\",\n  \"46-25=\": \"80-13=\",\n  \"87-75=\": \"26+38=\",\n  \"50-7=\": \"79-58=\",\n  \"20+51=\": \"14+58=\",\n  \"28-9=\": \"98-78=\",\n  \"78-39=\": \"32-15=\",\n  \"0+27=\": \"16+9=\",\n  \"97-76=\": \"26+70=\",\n  \"80+11=\": \"91-0=\",\n  \"61-5=\": \"22-16=\",\n  \"80-71=\": \"36+40=\",\n  \"14+1=\": \"5+6=\",\n  \"68-35=\": \"69-57=\",\n  \"60+9=\": \"76-69=\",\n  \"78-3=\": \"26+58=\",\n  \"23-19=\": \"84-55=\",\n  \"22+57=\": \"75+10=\",\n  \"48+10=\": \"50-35=\",\n  \"46-12=\": \"98-73=\",\n  \"83-69=\": \"30-20=\",\n  \"29+44=\": \"36+56=\",\n  \"55-18=\": \"90+6=\",\n  \"60-14=\": \"2+17=\",\n  \"96-6=\": \"92-35=\",\n  \"39+36=\": \"99-48=\",\n  \"98-47=\": \"98-81=\",\n  \"56-25=\": \"66+7=\",\n  \"96-34=\": \"52+16=\",\n  \"53+16=\": \"64+29=\",\n  \"61+19=\": \"47+30=\",\n  \"70-22=\": \"19+76=\",\n  \"75-6=\": \"65+4=\",\n  \"6+58=\": \"30+22=\",\n  \"66+24=\": \"58-11=\",\n  \"2+80=\": \"38+9=\",\n  \"48-29=\": \"28+25=\",\n  \"1+52=\": \"6+54=\",\n  \"5+32=\": \"64+35=\",\n  \"43+12=\": \"60-39=\",\n  \"26+61=\": \"88-31=\",\n  \"51-24=\": \"73-22=\",\n  \"18+44=\": \"6+10=\",\n  \"92-90=\": \"41+20=\",\n  \"18-14=\": \"69-22=\",\n  \"14+17=\": \"4+38=\",\n  \"50+26=\": \"55-46=\",\n  \"75-46=\": \"66-3=\",\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst oldValues = table.values;\nconst newValues = oldValues.map((row) =>\n  row.map((cellText) => {\n    const trimmed = cellText.trim();\n    return Object.prototype.hasOwnProperty.call(replacements, trimmed)\n      ? replacements[trimmed]\n      : cellText;\n  })\n);\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace every arithmetic-expression cell in the worksheet table with its\n# new expression. The mapping below is old-expression -> new-expression,\n# taken from the authoritative diff; using a lookup keyed on the existing\n# cell text (rather than blind positional indices) makes the edit robust to\n# the cells being visited in any order.\n$replacements = @{\n  \"44+39=\" = \"39+45=\";\n  \"65+10=\" = \"9-5=\";\n  \"56-27=\" = \"13-0=\";\n  \"77+0=\" = \"33+61=\";\n  \"13+9=\" = \"17+19=\";\n  \"88-63=\" = \"76-35=\";\n  \"59-6=\" = \"78-9=\";\n  \"63-36=\" = \"12+73=\";\n  \"21+54=\" = \"51-17=\";\n  \"26+73=\" = \"86-66=\";\n  \"52-37=\" = \"6+49=\";\n  \"37+19=\" = \"16-1=\";\n  \"84-15=\" = \"49-44=\";\n  \"83-68=\" = \"42+35=\";\n  \"5+33=\" = \"6+75=\";\n  \"97-5=\" = \"92-24=\";\n  \"11+44=\" = \"20-8=\";\n  \"30-26=\" = \"32+37=\";\n  \"66-52=\" = \"8+17=\";\n  \"1+67=\" = \"10+20=\";\n  \"22-17=\" = \"20-10=\";\n  \"37+20=\" = \"51-42=\";\n  \"4+88=\" = \"2+61=\";\n  \"26+42=\" = \"54+12=\";\n  \"57-47=\" = \"48-10=\";\n  \"14-7=\" = \"20+49=\";\n  \"11-0=\" = \"11+63=\";\n  \"82-5=\" = \"12+41=\";\n  \"91-66=\" = \"64+19=\";\n  \"66-44=\" = \"97-27=\";\n  \"36+31=\" = \"98-24=\";\n  \"83-67=\" = \"78-31=\";\n  \"66-4=\" = \"91-8=\";\n  \"65+11=\" = \"1+92=\";\n  \"37+32=\" = \"64-34=\";\n  \"51-5=\" = \"14+79=\";\n  \"91-52=\" = \"96-82=\";\n  \"33+37=\" = \"61-36=\";\n  \"57-7=\" = \"67-23=\";\n  \"1+61=\" = \"13+61=\";\n  \"32+60=\" = \"23+45=\";\n  \"39+43=\" = \"83+2=\";\n  \"22+58=\" = \"46-4=\";\n  \"94-6=\" = \"84-43=\";\n  \"37+25=\" = \"4+66=\";\n  \"64-20=\" = \"47-11=\";\n  \"32+4=\" = \"49+8=\";\n  \"13+69=\" = \"66-49=\";\n  \"45-8=\" = \"88-22=\";\n  \"72-52=\" = \"63-43=\";\n  \"29+8=\" = \"56+5=\";\n  \"38+31=\" = \"25-7=\";\n  \"75-25=\" = \"14+33=\";\n  \"46-25=\" = \"80-13=\";\n  \"87-75=\" = \"26+38=\";\n  \"50-7=\" = \"79-58=\";\n  \"20+51=\" = \"14+58=\";\n  \"28-9=\" = \"98-78=\";\n  \"78-39=\" = \"32-15=\";\n  \"0+27=\" = \"16+9=\";\n  \"97-76=\" = \"26+70=\";\n  \"80+11=\" = \"91-0=\";\n  \"61-5=\" = \"22-16=\";\n  \"80-71=\" = \"36+40=\";\n  \"14+1=\" = \"5+6=\";\n  \"68-35=\" = \"69-57=\";\n  \"60+9=\" = \"76-69=\";\n  \"78-3=\" = \"26+58=\";\n  \"23-19=\" = \"84-55=\";\n  \"22+57=\" = \"75+10=\";\n  \"48+10=\" = \"50-35=\";\n  \"46-12=\" = \"98-73=\";\n  \"83-69=\" = \"30-20=\";\n  \"29+44=\" = \"36+56=\";\n  \"55-18=\" = \"90+6=\";\n  \"60-14=\" = \"2+17=\";\n  \"96-6=\" = \"92-35=\";\n  \"39+36=\" = \"99-48=\";\n  \"98-47=\" = \"98-81=\";\n  \"56-25=\" = \"66+7=\";\n  \"96-34=\" = \"52+16=\";\n  \"53+16=\" = \"64+29=\";\n  \"61+19=\" = \"47+30=\";\n  \"70-22=\" = \"19+76=\";\n  \"75-6=\" = \"65+4=\";\n  \"6+58=\" = \"30+22=\";\n  \"66+24=\" = \"58-11=\";\n  \"2+80=\" = \"38+9=\";\n  \"48-29=\" = \"28+25=\";\n  \"1+52=\" = \"6+54=\";\n  \"5+32=\" = \"64+35=\";\n  \"43+12=\" = \"60-39=\";\n  \"26+61=\" = \"88-31=\";\n  \"51-24=\" = \"73-22=\";\n  \"18+44=\" = \"6+10=\";\n  \"92-90=\" = \"41+20=\";\n  \"18-14=\" = \"69-22=\";\n  \"14+17=\" = \"4+38=\";\n  \"50+26=\" = \"55-46=\";\n  \"75-46=\" = \"66-3=\";\n}\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cellText = $cell.Range.Text\n    # Strip the trailing cell-mark characters (CR + BEL) Word appends to\n    # Range.Text for table cells.\n    $key = $cellText.TrimEnd([char]13, [char]7)\n    if ($replacements.ContainsKey($key)) {\n      $cell.Range.Text = $replacements[$key]\n    }\n  }\n}\n"}
